$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing, so numeric-looking
# strings like "0.999" or "52.10" are preserved verbatim as text, matching
# the source inline-string cells, instead of Excel auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.313.10"
$ws.Range("E2").Value = "  +4.75%  "

$ws.Range("D3").Value = "2.246.04"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "229.11"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").Value = "0.636"
$ws.Range("E6").Value = "  +2.04%  "

$ws.Range("D7").Value = "65.19"
$ws.Range("E7").Value = "  +1.41%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +2.96%  "

$ws.Range("D10").Value = "0.0890"
$ws.Range("E10").Value = "  +3.38%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "2.580.32"
$ws.Range("E12").Value = "  +4.07%  "

$ws.Range("D13").Value = "16.19"
$ws.Range("E13").Value = "  +1.00%  "

$ws.Range("D14").Value = "22.42"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "0.831"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").Value = "2.242.14"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").Value = "41.225.53"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").Value = "74.22"
$ws.Range("E19").Value = "  +3.41%  "

$ws.Range("D20").Value = "0.0₃0914"
$ws.Range("E20").Value = "  +7.05%  "

$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22").Value = "257.13"
$ws.Range("E22").Value = "  +11.02%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -5.89%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").Value = "9.78"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").Value = "173.34"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").Value = "0.146"
$ws.Range("E28").Value = "  +4.31%  "

$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  +3.04%  "

$ws.Range("D30").Value = "1.45"
$ws.Range("E30").Value = "  +2.58%  "

$ws.Range("D31").Value = "2.83"
$ws.Range("E31").Value = "  +5.73%  "

$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").Value = "4.69"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").Value = "4.88"
$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").Value = "7.25"
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("D36").Value = "0.0634"
$ws.Range("E36").Value = "  +2.50%  "

$ws.Range("D37").Value = "3.84"
$ws.Range("E37").Value = "  +7.05%  "

$ws.Range("E38").Value = "  +3.36%  "

$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D40").Value = "0.000236"
$ws.Range("E40").Value = "  +55.75%  "

$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").Value = "4.89"
$ws.Range("E41").Value = "  +15.98%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0237"
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "8.80"
$ws.Range("E43").Value = "  +11.49%  "

$ws.Range("D44").Value = "102.36"
$ws.Range("E44").Value = "  -1.75%  "

$ws.Range("E45").Value = "  +5.25%  "

$ws.Range("D46").Value = "17.66"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "1.517.02"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").Value = "0.0942"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").Value = "52.10"
$ws.Range("E51").Value = "  +12.44%  "

# Restore column D to the default (General/Normal) style so the text-format
# override above does not leave a lingering NumberFormat on the cells.
$ws.Range("D2:D51").Style = "Normal"
